$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.121.73"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "3.063.19"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'527.44"
$ws.Range("E5").Value = "  +6.13%  "
$ws.Range("D6").Value = "'143.58"
$ws.Range("E6").Value = "  +6.81%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +5.49%  "
$ws.Range("D9").Value = "'7.63"
$ws.Range("E9").Value = "  +5.39%  "
$ws.Range("E10").Value = "  +8.02%  "
$ws.Range("E11").Value = "  +6.52%  "
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "3.589.14"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("E14").Value = "  +8.91%  "
$ws.Range("E15").Value = "  +17.33%  "
$ws.Range("D16").Value = "58.049.82"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("E17").Value = "  +9.13%  "
$ws.Range("D18").Value = "3.065.55"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  +7.29%  "
$ws.Range("D20").Value = "'8.21"
$ws.Range("E20").Value = "  +6.02%  "
$ws.Range("D21").Value = "'341.79"
$ws.Range("E21").Value = "  +5.04%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'5.74"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +8.28%  "
$ws.Range("D25").Value = "'65.49"
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("D26").Value = "0.0₃0976"
$ws.Range("E26").Value = "  +9.68%  "
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  +9.67%  "
$ws.Range("E30").Value = "  +10.65%  "
$ws.Range("E31").Value = "  +7.58%  "
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = "  +6.60%  "
$ws.Range("D33").Value = "'21.22"
$ws.Range("E33").Value = "  +4.67%  "
$ws.Range("E34").Value = "  +9.30%  "
$ws.Range("D35").Value = "'157.69"
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "  +7.47%  "
$ws.Range("D37").Value = "'1.33"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").Value = "'26.19"
$ws.Range("E38").Value = "  +13.96%  "
$ws.Range("D39").Value = "'0.0701"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").Value = "3.098.37"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("E42").Value = "  +12.90%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.669"
$ws.Range("E43").Value = "  +4.91%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.48"
$ws.Range("E45").Value = "  +6.64%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.05"
$ws.Range("E46").Value = "  +5.60%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.338.54"
$ws.Range("E47").Value = "  +5.13%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.01"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.10"
$ws.Range("E49").Value = "  +6.69%  "
$ws.Range("D50").Value = "'0.0244"
$ws.Range("E50").Value = "  +4.05%  "
$ws.Range("D51").Value = "'20.33"
$ws.Range("E51").Value = "  +7.32%  "

# Reset style on cells that were forced to text via quote-prefix,
# so no explicit style index is left on the cell (matches plain inline-string cells).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
